$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.04763786555579896
$ws.Range("C2").Value = 0.3127903958511391
$ws.Range("D2").Value = 0.1575252929769615
$ws.Range("E2").Value = 0.496779210170732
$ws.Range("F2").Value = 0
$ws.Range("G2").Value = 1.014732764554632
